$d = $word.ActiveDocument

# Locate the paragraph that currently holds "今天天气不错!心情也好！" -
# this is the last paragraph in the document (it also carries the
# _GoBack bookmark).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*今天天气不错!心情也好！*") {
        $targetIndex = $i
    }
}

# The paragraph right before the target (it already ends with an
# eastAsia-hinted paragraph mark, which is what the two new paragraphs
# need to inherit).
$prevIndex = $targetIndex - 1
$prevPara = $d.Paragraphs($prevIndex)

# Insert a new paragraph after it, add back the tab stop the target
# paragraph uses, and give it the "今天天气不错!心情也好！" text that
# used to live in the target paragraph.
$prevPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs($prevIndex + 1)
$newPara1.Range.ParagraphFormat.TabStops.Add(122.05)
$newPara1.Range.Text = "今天天气不错!心情也好！"

# Insert a second new paragraph after the first one (it inherits the
# tab stop already present on the paragraph it is split from).
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs($prevIndex + 2)
$newPara2.Range.Text = "哇哇哇！"

# Finally, replace the text of the original (now shifted) last
# paragraph with the new content, keeping its bookmark intact.
$finalIndex = $prevIndex + 3
$finalPara = $d.Paragraphs($finalIndex)
$finalPara.Range.Text = "今天学习了分支管理，创建了一个dev分支。使用Git创建分支简单又便捷"
